$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Apply 1.15 ("multiple") line spacing to every paragraph in the document.
#    w:line="276" w:lineRule="auto" == LineSpacingRule wdLineSpaceMultiple (5)
#    with LineSpacing expressed as 12pt * 1.15 = 13.8.
# ---------------------------------------------------------------------------
$d.Paragraphs.LineSpacingRule = 5
$d.Paragraphs.LineSpacing = 13.8

# ---------------------------------------------------------------------------
# 2. Underline the "Home Page Blocks:" heading paragraph.
# ---------------------------------------------------------------------------
$homeBlocksPara = $d.Paragraphs(5)
$homeBlocksPara.Range.Font.Underline = 1

# ---------------------------------------------------------------------------
# 3. Expand "Tai:" into the full bio line for Tai Shanahan, keeping "Tai
#    Shanahan: " bold and making the rest of the sentence non-bold.
# ---------------------------------------------------------------------------
$taiPara = $d.Paragraphs(11)
$taiStart = $taiPara.Range.Start

# Clear bold across the whole paragraph (including its end-of-paragraph
# mark) before editing, so the mark itself ends up non-bold like the rest
# of the new descriptive sentence; we re-apply bold to the name/colon
# segments afterwards.
$taiPara.Range.Font.Bold = 0
$taiPara = $d.Paragraphs(11)

# Replace the original "Tai:" text with the full sentence. The replaced
# range inherits the (now non-bold) character formatting; bold is restored
# per-segment below.
$oldTextRange = $d.Range($taiStart, $taiStart + 4)
$newText = "Tai Shanahan: Third-year history and political science major from Philadelphia, PA with concentrations in American politics and historical political violence."
$oldTextRange.Text = $newText

# Segment boundaries (relative offsets into $newText) and whether that
# segment should be bold. Each segment becomes its own <w:r> run.
$segments = @(
    @(0, 3, 1),     # "Tai"
    @(3, 12, 1),    # " Shanahan"
    @(12, 13, 1),   # ":"
    @(13, 14, 1),   # " "
    @(14, 82, 0),   # "Third-year history and political science major from Philadelphia, PA"
    @(82, 87, 0),   # " with"
    @(87, 88, 0),   # " "
    @(88, 158, 0)   # "concentrations in American politics and historical political violence."
)

foreach ($seg in $segments) {
    $segStart = $taiStart + $seg[0]
    $segEnd = $taiStart + $seg[1]
    $bold = $seg[2]
    $segRange = $d.Range($segStart, $segEnd)
    # Flip through the opposite state first so Word always records a fresh
    # run boundary here, even for the runs that stay bold.
    $segRange.Font.Bold = 1 - $bold
    $segRange.Font.Bold = $bold
}

Write-Host "done"
